$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells: replace spaces with underscores in subject names
$ws.Range("B1").Value = "General_Physics"
$ws.Range("C1").Value = "Classical_Mechanics"

# Move the active selection to I11 (as reflected in the saved view state)
$ws.Range("I11").Select()
